$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.739.98'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '3.482.62'
$ws.Range('E3').Value = '  -1.85%  '
$ws.Range('E4').Value = '  +0.11%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '603.12'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -2.37%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '147.66'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -4.11%  '
$ws.Range('D7').Value = '3.478.70'
$ws.Range('E7').Value = '  -1.82%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('E10').Value = '  -2.43%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '7.56'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  +3.46%  '
$ws.Range('E12').Value = '  -3.01%  '
$ws.Range('E13').Value = '  -3.06%  '
$ws.Range('D14').Value = '4.071.11'
$ws.Range('E14').Value = '  -1.86%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '31.43'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -4.58%  '
$ws.Range('D16').Value = '3.480.09'
$ws.Range('E16').Value = '  -1.85%  '
$ws.Range('D17').Value = '66.739.20'
$ws.Range('E17').Value = '  -1.65%  '
$ws.Range('E18').Value = '  +0.17%  '
$ws.Range('E19').Value = '  -4.58%  '
$ws.Range('E20').Value = '  -3.30%  '
$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '10.06'
$cell.Style = 'Normal'
$ws.Range('E21').Value = '  +1.18%  '
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '438.06'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  -3.37%  '
$ws.Range('E23').Value = '  -4.71%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '79.45'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '3.620.86'
$ws.Range('E26').Value = '  -1.77%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '0.0000120'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -7.52%  '
$ws.Range('E28').Value = '  -6.52%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '8.39'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -7.23%  '
$ws.Range('E30').Value = '  -2.90%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '1.58'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -5.54%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.168'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -0.92%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  +0.25%  '
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').Value = '3.473.20'
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('E36').Value = '  -5.54%  '
$ws.Range('E37').Value = '  -6.18%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '7.92'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -3.72%  '
$ws.Range('B39').Value = 'USDe'
$ws.Range('C39').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('E40').Value = '  +0.10%  '
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '177.11'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  -0.90%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.0884'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.35%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '2.12'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -9.87%  '
$ws.Range('E44').Value = '  -2.50%  '
$ws.Range('E45').Value = '  +0.03%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '46.37'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +1.74%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '28.82'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -6.67%  '
$ws.Range('E48').Value = '  -7.77%  '
$ws.Range('E49').Value = '  -4.17%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '2.43'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -7.92%  '
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '0.980'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -3.71%  '
